$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0
$batsman = "Shivam Mavi" + $nbsp

# Format the new rows as text first so the numeric-looking values
# (runs/balls/4s/6s/strike-rate) are stored as text, matching the rest
# of the sheet ("numberStoredAsText" ignored error covers A1:K7 afterwards).
$ws.Range("A5:K7").NumberFormat = "@"

# Row 5: duplicate of row 4 (Abu Dhabi / Sep 23 2020 vs Mumbai Indians)
$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " September 23 2020"
$ws.Range("C5").Value = "Mumbai won by 49 runs"
$ws.Range("D5").Value = "Kolkata Knight Riders"
$ws.Range("E5").Value = "Mumbai Indians"
$ws.Range("F5").Value = $batsman
$ws.Range("G5").Value = "9"
$ws.Range("H5").Value = "10"
$ws.Range("I5").Value = "1"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "90.00"

# Row 6: duplicate of row 2 (Sharjah / Oct 03 2020 vs Delhi Capitals)
$ws.Range("A6").Value = " Sharjah"
$ws.Range("B6").Value = " October 03 2020"
$ws.Range("C6").Value = "Capitals won by 18 runs"
$ws.Range("D6").Value = "Kolkata Knight Riders"
$ws.Range("E6").Value = "Delhi Capitals"
$ws.Range("F6").Value = $batsman
$ws.Range("G6").Value = "1"
$ws.Range("H6").Value = "3"
$ws.Range("I6").Value = "0"
$ws.Range("J6").Value = "0"
$ws.Range("K6").Value = "33.33"

# Row 7: duplicate of row 3 (Abu Dhabi / Oct 07 2020 vs Chennai Super Kings)
$ws.Range("A7").Value = " Abu Dhabi"
$ws.Range("B7").Value = " October 07 2020"
$ws.Range("C7").Value = "KKR won by 10 runs"
$ws.Range("D7").Value = "Kolkata Knight Riders"
$ws.Range("E7").Value = "Chennai Super Kings"
$ws.Range("F7").Value = $batsman
$ws.Range("G7").Value = "0"
$ws.Range("H7").Value = "1"
$ws.Range("I7").Value = "0"
$ws.Range("J7").Value = "0"
$ws.Range("K7").Value = "0.00"
